$d = $word.ActiveDocument

# Insert a new paragraph after the last paragraph ("PHP.ini Go to line no 916...")
# which inherits that paragraph's formatting (numPr/list, fonts, spacing, etc.)
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

# Set the text of the newly created paragraph
$p1 = $d.Paragraphs.Last
$p1.Range.Text = "PHP.ini Enable open_ssl port"

# Insert another paragraph after that one, again inheriting formatting
$p1 = $d.Paragraphs.Last
$p1.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last
$p2.Range.Text = "No Echo Statement in any of the Controller or Model"
